{"js": "// Update the date line and every \"AAA\u00d7B=\" problem cell in the practice\n// table. Old -> new text pairs, in document order, exactly mirroring the\n// unified diff.\nconst pairs = [\n  [\"2024-01-23 Tuesday\", \"2024-01-24 Wednesday\"],\n  [\"846\u00d73=\", \"687\u00d72=\"],\n  [\"901\u00d76=\", \"279\u00d77=\"],\n  [\"929\u00d78=\", \"661\u00d74=\"],\n  [\"936\u00d74=\", \"280\u00d74=\"],\n  [\"304\u00d75=\", \"400\u00d79=\"],\n  [\"922\u00d78=\", \"402\u00d79=\"],\n  [\"182\u00d76=\", \"497\u00d77=\"],\n  [\"990\u00d78=\", \"695\u00d77=\"],\n  [\"272\u00d72=\", \"431\u00d77=\"],\n  [\"251\u00d74=\", \"525\u00d78=\"],\n  [\"601\u00d74=\", \"142\u00d74=\"],\n  [\"497\u00d77=\", \"961\u00d74=\"],\n  [\"362\u00d73=\", \"842\u00d75=\"],\n  [\"573\u00d78=\", \"769\u00d72=\"],\n  [\"428\u00d79=\", \"178\u00d74=\"],\n  [\"980\u00d76=\", \"694\u00d75=\"],\n  [\"941\u00d73=\", \"684\u00d77=\"],\n  [\"231\u00d78=\", \"300\u00d75=\"],\n  [\"485\u00d78=\", \"698\u00d72=\"],\n  [\"935\u00d78=\", \"650\u00d74=\"],\n  [\"269\u00d77=\", \"621\u00d72=\"],\n  [\"616\u00d78=\", \"632\u00d77=\"],\n  [\"482\u00d76=\", \"875\u00d74=\"],\n  [\"971\u00d72=\", \"488\u00d74=\"],\n  [\"906\u00d75=\", \"600\u00d78=\"],\n];\n\nconst body = context.document.body;\n\n// Resolve every \"old\" needle against the document's ORIGINAL text in one\n// batch, before any replacement is written. Several of the new values\n// collide with other entries' old values (e.g. one cell goes 182\u00d76= ->\n// 497\u00d77= while another goes 497\u00d77= -> 961\u00d74=), so searching must happen\n// up front -- otherwise a later search could accidentally match text a\n// prior step just inserted.\nconst searches = pairs.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: true })\n);\nsearches.forEach((s) => s.load(\"items\"));\nawait context.sync();\n\nfor (let i = 0; i < pairs.length; i++) {\n  const [oldText, newText] = pairs[i];\n  const items = searches[i].items;\n  if (items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n  items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every \"AAA\u00d7B=\" problem cell in the practice\n# table, mirroring the unified diff exactly.\n#\n# NOTE: several of the new values collide with other cells' old values\n# (e.g. row5/col2 goes 182\u00d76= -> 497\u00d77= while row10/col2 goes\n# 497\u00d77= -> 961\u00d74=), so a blind Find/Replace-All pass across the whole\n# document would clobber both cells once the first replacement lands.\n# Addressing each cell by its (row, col) table position sidesteps that\n# entirely -- no text search is involved.\n\n$d = $word.ActiveDocument\n\n# --- Date line (first paragraph, outside the table) ---\n$dateRange = $d.Paragraphs(1).Range\n$dateRange.MoveEnd(1, -1) | Out-Null   # drop the trailing paragraph mark\nif ($dateRange.Text -ne \"2024-01-23 Tuesday\") {\n    throw \"Unexpected date text: '$($dateRange.Text)'\"\n}\n$dateRange.Text = \"2024-01-24 Wednesday\"\n\n# --- Table cells (5 columns; data lives in rows 1, 5, 10, 15, 20) ---\n$cells = @(\n  @{ Row = 1; Col = 1; Old = \"846\u00d73=\"; New = \"687\u00d72=\" },\n  @{ Row = 1; Col = 2; Old = \"901\u00d76=\"; New = \"279\u00d77=\" },\n  @{ Row = 1; Col = 3; Old = \"929\u00d78=\"; New = \"661\u00d74=\" },\n  @{ Row = 1; Col = 4; Old = \"936\u00d74=\"; New = \"280\u00d74=\" },\n  @{ Row = 1; Col = 5; Old = \"304\u00d75=\"; New = \"400\u00d79=\" },\n  @{ Row = 5; Col = 1; Old = \"922\u00d78=\"; New = \"402\u00d79=\" },\n  @{ Row = 5; Col = 2; Old = \"182\u00d76=\"; New = \"497\u00d77=\" },\n  @{ Row = 5; Col = 3; Old = \"990\u00d78=\"; New = \"695\u00d77=\" },\n  @{ Row = 5; Col = 4; Old = \"272\u00d72=\"; New = \"431\u00d77=\" },\n  @{ Row = 5; Col = 5; Old = \"251\u00d74=\"; New = \"525\u00d78=\" },\n  @{ Row = 10; Col = 1; Old = \"601\u00d74=\"; New = \"142\u00d74=\" },\n  @{ Row = 10; Col = 2; Old = \"497\u00d77=\"; New = \"961\u00d74=\" },\n  @{ Row = 10; Col = 3; Old = \"362\u00d73=\"; New = \"842\u00d75=\" },\n  @{ Row = 10; Col = 4; Old = \"573\u00d78=\"; New = \"769\u00d72=\" },\n  @{ Row = 10; Col = 5; Old = \"428\u00d79=\"; New = \"178\u00d74=\" },\n  @{ Row = 15; Col = 1; Old = \"980\u00d76=\"; New = \"694\u00d75=\" },\n  @{ Row = 15; Col = 2; Old = \"941\u00d73=\"; New = \"684\u00d77=\" },\n  @{ Row = 15; Col = 3; Old = \"231\u00d78=\"; New = \"300\u00d75=\" },\n  @{ Row = 15; Col = 4; Old = \"485\u00d78=\"; New = \"698\u00d72=\" },\n  @{ Row = 15; Col = 5; Old = \"935\u00d78=\"; New = \"650\u00d74=\" },\n  @{ Row = 20; Col = 1; Old = \"269\u00d77=\"; New = \"621\u00d72=\" },\n  @{ Row = 20; Col = 2; Old = \"616\u00d78=\"; New = \"632\u00d77=\" },\n  @{ Row = 20; Col = 3; Old = \"482\u00d76=\"; New = \"875\u00d74=\" },\n  @{ Row = 20; Col = 4; Old = \"971\u00d72=\"; New = \"488\u00d74=\" },\n  @{ Row = 20; Col = 5; Old = \"906\u00d75=\"; New = \"600\u00d78=\" }\n)\n\n$table = $d.Tables(1)\n\nforeach ($item in $cells) {\n    $cell = $table.Cell($item.Row, $item.Col)\n    $rng = $cell.Range\n    $rng.MoveEnd(1, -1) | Out-Null   # drop the cell-end mark\n    if ($rng.Text -ne $item.Old) {\n        throw \"Cell ($($item.Row),$($item.Col)) expected '$($item.Old)' but found '$($rng.Text)'\"\n    }\n    $rng.Text = $item.New\n}\n"}
